$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value2 = "'67.822.75"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.12%  "
$c = $ws.Range("D3")
$c.Value2 = "'3.500.26"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.13%  "
$c = $ws.Range("D4")
$c.Value2 = "'1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$c = $ws.Range("D5")
$c.Value2 = "'606.68"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.32%  "
$c = $ws.Range("D6")
$c.Value2 = "'150.81"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.33%  "
$c = $ws.Range("D7")
$c.Value2 = "'3.497.39"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("E8").Value = "  +0.03%  "
$c = $ws.Range("D9")
$c.Value2 = "'0.487"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.45%  "
$c = $ws.Range("D10")
$c.Value2 = "'0.144"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.44%  "
$c = $ws.Range("D11")
$c.Value2 = "'7.57"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +6.98%  "
$c = $ws.Range("D12")
$c.Value2 = "'0.431"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.91%  "
$c = $ws.Range("D13")
$c.Value2 = "'0.0000215"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.19%  "
$c = $ws.Range("D14")
$c.Value2 = "'32.04"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.25%  "
$c = $ws.Range("D15")
$c.Value2 = "'4.095.54"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.06%  "
$c = $ws.Range("D16")
$c.Value2 = "'3.507.21"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.07%  "
$c = $ws.Range("D17")
$c.Value2 = "'67.853.71"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("E18").Value = "  -0.24%  "
$c = $ws.Range("D19")
$c.Value2 = "'6.49"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.09%  "
$c = $ws.Range("D20")
$c.Value2 = "'15.38"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.05%  "
$c = $ws.Range("D21")
$c.Value2 = "'9.93"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.13%  "
$c = $ws.Range("D22")
$c.Value2 = "'446.01"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.44%  "
$c = $ws.Range("D23")
$c.Value2 = "'0.625"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "
$c = $ws.Range("D24")
$c.Value2 = "'79.18"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.13%  "
$c = $ws.Range("D25")
$c.Value2 = "'3.640.87"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D26")
$c.Value2 = "'1.00"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c = $ws.Range("D27")
$c.Value2 = "'0.0000126"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.95%  "
$c = $ws.Range("D28")
$c.Value2 = "'8.66"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.20%  "
$c = $ws.Range("D29")
$c.Value2 = "'9.95"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -3.10%  "
$ws.Range("E30").Value = "  -1.49%  "
$c = $ws.Range("D31")
$c.Value2 = "'1.64"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("E33").Value = "  -0.06%  "
$c = $ws.Range("D34")
$c.Value2 = "'25.60"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.47%  "
$c = $ws.Range("D35")
$c.Value2 = "'6.14"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.51%  "
$c = $ws.Range("D36")
$c.Value2 = "'1.85"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.49%  "
$c = $ws.Range("D37")
$c.Value2 = "'3.497.10"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.87%  "
$c = $ws.Range("D38")
$c.Value2 = "'7.99"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("E39").Value = "  -0.02%  "
$c = $ws.Range("D40")
$c.Value2 = "'2.33"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +5.28%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D41")
$c.Value2 = "'1.00"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D42")
$c.Value2 = "'177.08"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.88%  "
$c = $ws.Range("D43")
$c.Value2 = "'0.0897"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "
$c = $ws.Range("D44")
$c.Value2 = "'5.42"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D45")
$c.Value2 = "'0.893"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D46")
$c.Value2 = "'30.26"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +4.71%  "
$c = $ws.Range("D47")
$c.Value2 = "'46.59"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("E48").Value = "  +0.38%  "
$c = $ws.Range("D49")
$c.Value2 = "'2.53"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -5.42%  "
$c = $ws.Range("D50")
$c.Value2 = "'7.61"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.31%  "
$c = $ws.Range("D51")
$c.Value2 = "'0.989"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.55%  "
